# "Ispravljen dizajn izvjesca (ne prekida stranicu usred tablice)"
#
# The "Kupci" sheet table is restructured so it no longer breaks awkwardly:
# the old leading index column ("oznaka_poreznog_broja", values 1/1) is
# dropped, the tax-number column becomes the new first column, two new
# (currently empty-bodied) columns are inserted for "pdv_identifikacijski_broj"
# and "ostali_brojevi", and the customer-name column moves to the end.

$wb = $excel.ActiveWorkbook

$wsZaglavlje = $wb.Worksheets.Item("Zaglavlje")
$wsKupci     = $wb.Worksheets.Item("Kupci")
$wsRacuni    = $wb.Worksheets.Item("Racuni")

# --- Restructure the "Kupci" sheet -----------------------------------------

# Drop the old first column (header "oznaka_poreznog_broja", values 1/1) -
# it shifts "porezni_broj" into column A and "naziv_kupca" into column B.
$wsKupci.Columns.Item(1).Delete()

# Make room for the two new columns between "porezni_broj" (A) and
# "naziv_kupca" (now B, destined for D).
$wsKupci.Columns.Item(2).Insert()
$wsKupci.Columns.Item(3).Insert()

# The inserted columns picked up neighbouring formatting on the body rows -
# the new columns have no data rows, only a header.
$wsKupci.Range("B2:C3").Clear()

$wsKupci.Range("B1").ClearFormats()
$wsKupci.Range("B1").Font.Bold = $true
$wsKupci.Range("B1").Value = "pdv_identifikacijski_broj"

$wsKupci.Range("C1").ClearFormats()
$wsKupci.Range("C1").Font.Bold = $true
$wsKupci.Range("C1").Value = "ostali_brojevi"

# Column widths for the two new columns (best-fit equivalents).
$wsKupci.Columns.Item(2).ColumnWidth = 20.666666666666668  # -> stored width 21.5
$wsKupci.Columns.Item(3).ColumnWidth = 11.666666666666666  # -> stored width 12.5

# --- View / selection state --------------------------------------------

# Window moved slightly and the active tab switches from "Zaglavlje" to
# "Kupci".
$win = $wb.Windows.Item(1)
$win.Left = 780

$wsZaglavlje.Range("D22").Select()

$wsKupci.Activate()
$wsKupci.Range("E12").Select()

# "Racuni" keeps its previous selection untouched.
